$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.026.63'
$ws.Range("E2").Value = '  -0.16%  '

# Row 3
$ws.Range("D3").Value = '2.367.28'
$ws.Range("E3").Value = '  -3.47%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.39'
$ws.Range("E5").Value = '  -0.86%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.28'
$ws.Range("E6").Value = '  -5.20%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("E8").Value = '  -10.38%  '

# Row 9
$ws.Range("D9").Value = '2.367.26'
$ws.Range("E9").Value = '  -3.50%  '

# Row 10
$ws.Range("E10").Value = '  -0.95%  '

# Row 11
$ws.Range("E11").Value = '  +0.36%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.24'
$ws.Range("E12").Value = '  -3.02%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.343'
$ws.Range("E13").Value = '  -2.17%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.98'
$ws.Range("E14").Value = '  -3.45%  '

# Row 15
$ws.Range("D15").Value = '2.790.83'
$ws.Range("E15").Value = '  -3.49%  '

# Row 16
$ws.Range("E16").Value = '  -1.55%  '

# Row 17
$ws.Range("D17").Value = '60.862.51'
$ws.Range("E17").Value = '  -0.33%  '

# Row 18
$ws.Range("D18").Value = '2.362.88'
$ws.Range("E18").Value = '  -3.57%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.64'
$ws.Range("E19").Value = '  -3.62%  '

# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '316.19'
$ws.Range("E20").Value = '  +0.01%  '

# Row 21
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.07'
$ws.Range("E21").Value = '  -1.67%  '

# Row 22
$ws.Range("E22").Value = '  -5.93%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.14%  '

# Row 24
$ws.Range("E24").Value = '  +2.92%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.42'
$ws.Range("E25").Value = '  -0.33%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.50'
$ws.Range("E26").Value = '  +13.00%  '

# Row 27
$ws.Range("E27").Value = '  -0.09%  '

# Row 28
$ws.Range("D28").Value = '2.483.63'
$ws.Range("E28").Value = '  -3.71%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0906'
$ws.Range("E29").Value = '  -5.39%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.99'
$ws.Range("E30").Value = '  -2.28%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '507.39'
$ws.Range("E31").Value = '  -7.14%  '

# Row 32
$ws.Range("E32").Value = '  -3.78%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.146'
$ws.Range("E33").Value = '  +0.38%  '

# Row 34
$ws.Range("E34").Value = '  -5.01%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("E35").Value = '  -1.91%  '

# Row 36
$ws.Range("E36").Value = '  +0.13%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.59'
$ws.Range("E37").Value = '  -3.99%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.57'
$ws.Range("E38").Value = '  +1.39%  '

# Row 39
$ws.Range("E39").Value = '  -1.16%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.28'
$ws.Range("E40").Value = '  -9.15%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.78'
$ws.Range("E41").Value = '  +1.94%  '

# Row 42
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.47'
$ws.Range("E42").Value = '  -1.63%  '

# Row 43
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.04%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.18'
$ws.Range("E44").Value = '  -0.30%  '

# Row 45
$ws.Range("E45").Value = '  -7.16%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '138.97'
$ws.Range("E46").Value = '  -4.67%  '

# Row 47
$ws.Range("E47").Value = '  -1.05%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0514'
$ws.Range("E48").Value = '  -3.40%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.60'
$ws.Range("E49").Value = '  -8.01%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.573'
$ws.Range("E50").Value = '  -2.11%  '

# Row 51
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0898'
$ws.Range("E51").Value = '  -3.44%  '
